# Applies the "Add examples and validation for first assay type" change:
#  - Instructions sheet: bold the title cell (A1) and fix an apostrophe typo (A6)
#  - Antibodies sheet: extend the two list-validation ranges by one row
#  - Terminology sheet: add a Host example (A4), rename sIgA->kappa (B16),
#    and add a new "lambda" row (17)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Instructions sheet
# ---------------------------------------------------------------------
$instructions = $wb.Worksheets.Item("Instructions")
$instructions.Unprotect()

# Make the title bold by copying the existing bold header style from the
# Antibodies sheet (keeps the same shared cell style instead of minting one).
$antibodies = $wb.Worksheets.Item("Antibodies")
$antibodies.Range("A1").Copy()
$instructions.Range("A1").PasteSpecial(-4122)

# Fix "Your institutions" -> "Your institution's"
$instructions.Range("A6").Value = "- Antibody name: Your institution's preferred name for the antibody."

$instructions.Protect()

# ---------------------------------------------------------------------
# 2) Antibodies sheet: widen the two data validation source ranges
# ---------------------------------------------------------------------
$antibodies.Range("B2:B100").Validation.Modify(3, 1, 1, "=Terminology!A2:A4")
$antibodies.Range("C2:C100").Validation.Modify(3, 1, 1, "=Terminology!B2:B17")

# ---------------------------------------------------------------------
# 3) Terminology sheet: new host example, rename a value, add a new row
# ---------------------------------------------------------------------
$terminology = $wb.Worksheets.Item("Terminology")
$terminology.Unprotect()

$terminology.Range("A4").Value = "Mus musculus BALB/C"
$terminology.Range("B16").Value = "kappa"

# Row 17: A17 stays blank (matching the blank A-column cells above it),
# B17 gets the new "lambda" value.
$terminology.Range("A17").Style = "Normal"
$terminology.Range("B17").Value = "lambda"

$terminology.Protect()
